# Auto-generated Excel COM-interop script applying scheduled price/profit refresh
# to the Masamune_Profits workbook (per commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 60.875
$ws.Range("I2").Value = 60.875
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 60.875
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 52.125
$ws.Range("N2").ClearContents()

# ALC row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 233336670
$ws.Range("I29").Value = 100005000
$ws.Range("J29").Value = 500000000
$ws.Range("K29").Value = 300015000
$ws.Range("L29").Value = 1500000000
$ws.Range("M29").Value = -300014719
$ws.Range("N29").Value = -1500000562

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 71428670
$ws.Range("I38").Value = 71428670
$ws.Range("K38").Value = 214286010
$ws.Range("M38").Value = -214285638

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 4762354.5
$ws.Range("I58").Value = 5952560
$ws.Range("K58").Value = 17857680
$ws.Range("M58").Value = -17857530

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27496

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -87480

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2142.8572
$ws.Range("I113").Value = 2160
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 2160
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 1094
$ws.Range("N113").Value = -8608

# ALC row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2384.08
$ws.Range("J131").Value = 3350.75
$ws.Range("L131").Value = 10052.25
$ws.Range("N131").Value = -20132.25

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 22986.512
$ws.Range("I132").Value = 2863.1538
$ws.Range("K132").Value = 8589.4614
$ws.Range("M132").Value = -6059.4614

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2202.028
$ws.Range("I138").Value = 1980.2222
$ws.Range("J138").Value = 2277.3584
$ws.Range("K138").Value = 5940.6666
$ws.Range("L138").Value = 6832.0752
$ws.Range("M138").Value = -800.6665999999996
$ws.Range("N138").Value = -17112.0752

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4431.2085
$ws.Range("I32").Value = 2087.1428
$ws.Range("J32").Value = 20839.666
$ws.Range("K32").Value = 2087.1428
$ws.Range("L32").Value = 20839.666
$ws.Range("M32").Value = -1800.1428
$ws.Range("N32").Value = -21413.666

# BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 9800
$ws.Range("J81").Value = 9800
$ws.Range("L81").Value = 9800
$ws.Range("N81").Value = -11922

# BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 9800
$ws.Range("J84").Value = 9800
$ws.Range("L84").Value = 29400
$ws.Range("N84").Value = -40008

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3286.4348
$ws.Range("I99").Value = 3054.75
$ws.Range("J99").Value = 3410
$ws.Range("K99").Value = 3054.75
$ws.Range("L99").Value = 3410
$ws.Range("M99").Value = -1556.75
$ws.Range("N99").Value = -6406

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3286.4348
$ws.Range("I126").Value = 3054.75
$ws.Range("J126").Value = 3410
$ws.Range("K126").Value = 9164.25
$ws.Range("L126").Value = 10230
$ws.Range("M126").Value = -6694.25
$ws.Range("N126").Value = -15170

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6919.9375
$ws.Range("I5").Value = 7247.933
$ws.Range("K5").Value = 21743.799
$ws.Range("M5").Value = -21631.799

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 200303.5
$ws.Range("I38").Value = 286.25
$ws.Range("J38").Value = 333648.34
$ws.Range("K38").Value = 858.75
$ws.Range("L38").Value = 1000945.02
$ws.Range("M38").Value = -511.75
$ws.Range("N38").Value = -1001639.02

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5833.8335
$ws.Range("I68").Value = 1000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2189

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 5833.8335
$ws.Range("I71").Value = 1000
$ws.Range("K71").Value = 9000
$ws.Range("M71").Value = -4944

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 6919.9375
$ws.Range("I135").Value = 7247.933
$ws.Range("K135").Value = 65231.397
$ws.Range("M135").Value = -62696.397

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6086.273
$ws.Range("I80").Value = 4778.5
$ws.Range("K80").Value = 4778.5
$ws.Range("M80").Value = -3780.5

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 6086.273
$ws.Range("I83").Value = 4778.5
$ws.Range("K83").Value = 23892.5
$ws.Range("M83").Value = -18900.5

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4049.3572
$ws.Range("I40").Value = 3350.1
$ws.Range("J40").Value = 5797.5
$ws.Range("K40").Value = 3350.1
$ws.Range("L40").Value = 5797.5
$ws.Range("M40").Value = -3214.1
$ws.Range("N40").Value = -6069.5

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 33334660
$ws.Range("J54").Value = 1990
$ws.Range("L54").Value = 1990
$ws.Range("N54").Value = -3030

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2990.1765
$ws.Range("I62").Value = 2958.25
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2958.25
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2334.25
$ws.Range("N62").Value = -4248

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 2990.1765
$ws.Range("I65").Value = 2958.25
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 14791.25
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11671.25
$ws.Range("N65").Value = -21240

# WVR row 76
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 36173
$ws.Range("J76").Value = 36173
$ws.Range("L76").Value = 36173
$ws.Range("N76").Value = -36803

# WVR row 79
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H79").Value = 36173
$ws.Range("J79").Value = 36173
$ws.Range("L79").Value = 36173
$ws.Range("N79").Value = -38357

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1398.5714
$ws.Range("I81").Value = 1398.5714
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2797.1428
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1736.1428
$ws.Range("N81").ClearContents()

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1398.5714
$ws.Range("I84").Value = 1398.5714
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 13985.714
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -8681.714
$ws.Range("N84").ClearContents()

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1403.1333
$ws.Range("I126").Value = 971.4091
$ws.Range("K126").Value = 2914.2273
$ws.Range("M126").Value = -444.2273
